$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 40.9
$ws.Range("I8").Value = 40.9
$ws.Range("K8").Value = 122.7
$ws.Range("M8").Value = 16.30000000000001

$ws.Range("H21").Value = 11431
$ws.Range("I21").Value = 1672.3334
$ws.Range("K21").Value = 1672.3334
$ws.Range("M21").Value = -1204.3334

$ws.Range("H23").Value = 11431
$ws.Range("I23").Value = 1672.3334
$ws.Range("K23").Value = 1672.3334
$ws.Range("M23").Value = -1438.3334

$ws.Range("H76").Value = 2587263.2
$ws.Range("I76").Value = 3587376.2
$ws.Range("J76").Value = 3637.3333
$ws.Range("K76").Value = 3587376.2
$ws.Range("L76").Value = 3637.3333
$ws.Range("M76").Value = -3587061.2
$ws.Range("N76").Value = -4267.3333

$ws.Range("H79").Value = 2587263.2
$ws.Range("I79").Value = 3587376.2
$ws.Range("J79").Value = 3637.3333
$ws.Range("K79").Value = 3587376.2
$ws.Range("L79").Value = 3637.3333
$ws.Range("M79").Value = -3586284.2
$ws.Range("N79").Value = -5821.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20057.525
$ws.Range("I32").Value = 2643.4546
$ws.Range("K32").Value = 2643.4546
$ws.Range("M32").Value = -2356.4546

$ws.Range("H45").Value = 1693.5333
$ws.Range("I45").Value = 1010.1667
$ws.Range("J45").Value = 2149.111
$ws.Range("K45").Value = 1010.1667
$ws.Range("L45").Value = 2149.111
$ws.Range("M45").Value = -633.1667
$ws.Range("N45").Value = -2903.111

$ws.Range("H104").Value = 29000
$ws.Range("J104").Value = 29000
$ws.Range("L104").Value = 29000
$ws.Range("N104").Value = -35988

$ws.Range("H132").Value = 3037.8936
$ws.Range("I132").Value = 2629.2896
$ws.Range("J132").Value = 4763.1113
$ws.Range("K132").Value = 7887.8688
$ws.Range("L132").Value = 14289.3339
$ws.Range("M132").Value = -5357.8688
$ws.Range("N132").Value = -19349.3339

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 33360
$ws.Range("I82").Value = 8500
$ws.Range("J82").Value = 49933.332
$ws.Range("K82").Value = 8500
$ws.Range("L82").Value = 49933.332
$ws.Range("M82").Value = -8117
$ws.Range("N82").Value = -50699.332

$ws.Range("H85").Value = 33360
$ws.Range("I85").Value = 8500
$ws.Range("J85").Value = 49933.332
$ws.Range("K85").Value = 8500
$ws.Range("L85").Value = 49933.332
$ws.Range("M85").Value = -7174
$ws.Range("N85").Value = -52585.332

$ws.Range("H94").Value = 1638.1538
$ws.Range("I94").Value = 1061.6
$ws.Range("J94").Value = 3560
$ws.Range("K94").Value = 1061.6
$ws.Range("L94").Value = 3560
$ws.Range("M94").Value = -610.5999999999999
$ws.Range("N94").Value = -4462

$ws.Range("H97").Value = 265449.5
$ws.Range("J97").Value = 30471
$ws.Range("L97").Value = 30471
$ws.Range("N97").Value = -32453

$ws.Range("H106").Value = 21654.2
$ws.Range("J106").Value = 21654.2
$ws.Range("L106").Value = 21654.2
$ws.Range("N106").Value = -24178.2

$ws.Range("H107").Value = 1717.2307
$ws.Range("I107").Value = 1437.25
$ws.Range("J107").Value = 2165.2
$ws.Range("K107").Value = 1437.25
$ws.Range("L107").Value = 2165.2
$ws.Range("M107").Value = 482.75
$ws.Range("N107").Value = -6005.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 258750
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 258750
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 258750
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -259028

$ws.Range("H23").Value = 40801.6
$ws.Range("I23").Value = 43003
$ws.Range("J23").Value = 37499.5
$ws.Range("K23").Value = 43003
$ws.Range("L23").Value = 37499.5
$ws.Range("M23").Value = -42763
$ws.Range("N23").Value = -37979.5

$ws.Range("H27").Value = 40801.6
$ws.Range("I27").Value = 43003
$ws.Range("J27").Value = 37499.5
$ws.Range("K27").Value = 43003
$ws.Range("L27").Value = 37499.5
$ws.Range("M27").Value = -42811
$ws.Range("N27").Value = -37883.5

$ws.Range("H31").Value = 3833.0889
$ws.Range("I31").Value = 1028.2703
$ws.Range("K31").Value = 1028.2703
$ws.Range("M31").Value = -733.2702999999999

$ws.Range("H34").Value = 3833.0889
$ws.Range("I34").Value = 1028.2703
$ws.Range("K34").Value = 1028.2703
$ws.Range("M34").Value = -826.2702999999999

$ws.Range("H107").Value = 582.5
$ws.Range("I107").Value = 298.33334
$ws.Range("J107").Value = 779.2308
$ws.Range("K107").Value = 298.33334
$ws.Range("L107").Value = 779.2308
$ws.Range("M107").Value = 1621.66666
$ws.Range("N107").Value = -4619.2308

$ws.Range("H134").Value = 3939.15
$ws.Range("I134").Value = 1973.1818
$ws.Range("K134").Value = 5919.5454
$ws.Range("M134").Value = -3384.5454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1300.381
$ws.Range("I51").Value = 2436
$ws.Range("J51").Value = 1111.1111
$ws.Range("K51").Value = 7308
$ws.Range("L51").Value = 3333.3333
$ws.Range("M51").Value = -6848
$ws.Range("N51").Value = -4253.3333

$ws.Range("H60").Value = 2135.5
$ws.Range("I60").Value = 215
$ws.Range("J60").Value = 2615.625
$ws.Range("K60").Value = 645
$ws.Range("L60").Value = 7846.875
$ws.Range("M60").Value = -394
$ws.Range("N60").Value = -8348.875

$ws.Range("H68").Value = 15998.375
$ws.Range("I68").Value = 20608.5
$ws.Range("J68").Value = 2168
$ws.Range("K68").Value = 61825.5
$ws.Range("L68").Value = 6504
$ws.Range("M68").Value = -61014.5
$ws.Range("N68").Value = -8126

$ws.Range("H71").Value = 15998.375
$ws.Range("I71").Value = 20608.5
$ws.Range("J71").Value = 2168
$ws.Range("K71").Value = 185476.5
$ws.Range("L71").Value = 19512
$ws.Range("M71").Value = -181420.5
$ws.Range("N71").Value = -27624

$ws.Range("H131").Value = 6290707
$ws.Range("J131").Value = 6411676.5
$ws.Range("L131").Value = 19235029.5
$ws.Range("N131").Value = -19245109.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4938.909
$ws.Range("I70").Value = 4717.273
$ws.Range("J70").Value = 5825.4546
$ws.Range("K70").Value = 4717.273
$ws.Range("L70").Value = 5825.4546
$ws.Range("M70").Value = -4447.273
$ws.Range("N70").Value = -6365.4546

$ws.Range("H73").Value = 4938.909
$ws.Range("I73").Value = 4717.273
$ws.Range("J73").Value = 5825.4546
$ws.Range("K73").Value = 4717.273
$ws.Range("L73").Value = 5825.4546
$ws.Range("M73").Value = -3781.273
$ws.Range("N73").Value = -7697.4546

$ws.Range("H80").Value = 2565.1853
$ws.Range("I80").Value = 2395.2
$ws.Range("J80").Value = 3050.8572
$ws.Range("K80").Value = 2395.2
$ws.Range("L80").Value = 3050.8572
$ws.Range("M80").Value = -1397.2
$ws.Range("N80").Value = -5046.8572

$ws.Range("H83").Value = 2565.1853
$ws.Range("I83").Value = 2395.2
$ws.Range("J83").Value = 3050.8572
$ws.Range("K83").Value = 11976
$ws.Range("L83").Value = 15254.286
$ws.Range("M83").Value = -6984
$ws.Range("N83").Value = -25238.286

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 923.5454999999999
$ws.Range("I107").Value = 923.5454999999999
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 923.5454999999999
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 996.4545000000001
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -830
$ws.Range("N17").Value = -3340

$ws.Range("H40").Value = 3275.8845
$ws.Range("I40").Value = 2990.6
$ws.Range("K40").Value = 2990.6
$ws.Range("M40").Value = -2854.6

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H106").Value = 23206.445
$ws.Range("J106").Value = 23206.445
$ws.Range("L106").Value = 23206.445
$ws.Range("N106").Value = -25730.445

$ws.Range("H122").Value = 2951.6667
$ws.Range("I122").Value = 2013.9286
$ws.Range("K122").Value = 6041.7858
$ws.Range("M122").Value = -3591.7858

$ws.Range("H132").Value = 3366.5
$ws.Range("I132").Value = 2380.85
$ws.Range("J132").Value = 5830.625
$ws.Range("K132").Value = 7142.549999999999
$ws.Range("L132").Value = 17491.875
$ws.Range("M132").Value = -4612.549999999999
$ws.Range("N132").Value = -22551.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 2623.75
$ws.Range("I17").Value = 2623.75
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2623.75
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -2451.75
$ws.Range("N17").ClearContents()

$ws.Range("H122").Value = 143986.28
$ws.Range("I122").Value = 200780.8
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 602342.3999999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -599892.3999999999
$ws.Range("N122").Value = -10900
